# The sheet had a stray "index" row (row 2: literal values 1..12, using a
# bold 8pt font that exists nowhere else in the sheet) sitting between the
# header row and the real data-entry row. Remove that whole row so the data
# row (and the totals row below it) shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()

# Restore the user's selection to the cell they ended up on after the edit.
$ws.Range("C19").Select()
